$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Valor Mora" / "Salario Basico" figures for the two most
# recently added workers (rows 17 = RONALD JAVIER PALENCIA LECHUGA and
# 18 = CRISTIAN ENRIQUE MORALES ORTIZ) as part of the EC database refresh.
$ws.Range("F17").Value = 489858
$ws.Range("G17").Value = 12246443
$ws.Range("F18").Value = 1333
$ws.Range("G18").Value = 1133000
